$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so that
# numeric-looking strings (e.g. "1.00", "624.79") are not silently
# converted to numbers by Excel's auto-detection, matching the original
# inline-string cell type used throughout this sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "90.562.20"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.126.69"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "234.38"
$ws.Range("E5").Value = "  +6.38%  "
$ws.Range("D6").Value = "624.79"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +11.19%  "
$ws.Range("E8").Value = "  -5.48%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.125.49"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").Value = "0.719"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("D13").Value = "36.17"
$ws.Range("E13").Value = "  +4.78%  "
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "5.60"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "90.169.57"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "3.686.45"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "3.068.07"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "3.67"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  -6.28%  "
$ws.Range("D22").Value = "449.56"
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("D23").Value = "9.02"
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").Value = "5.87"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "90.71"
$ws.Range("E26").Value = "  +8.05%  "
$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "3.259.74"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "9.25"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").Value = "0.159"
$ws.Range("E31").Value = "  -5.62%  "
$ws.Range("D32").Value = "27.75"
$ws.Range("E32").Value = "  +18.30%  "
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").Value = "0.195"
$ws.Range("E34").Value = "  +30.85%  "
$ws.Range("D35").Value = "0.152"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("D36").Value = "506.26"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("D37").Value = "3.67"
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "0.425"
$ws.Range("E41").Value = "  +11.93%  "
$ws.Range("D42").Value = "22.19"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +9.75%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "3.20"
$ws.Range("E45").Value = "  +29.52%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.97"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").Value = "0.699"
$ws.Range("E47").Value = "  +12.31%  "
$ws.Range("D48").Value = "148.96"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +9.44%  "
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").Value = "44.61"
$ws.Range("E51").Value = "  +0.85%  "

# Restore the default (Normal) style so no stray number-format
# attribute is left behind on these cells.
$ws.Range("D2:E51").Style = "Normal"
